# Atualizado por script em 01-12-2023 14:45
# Adds two new match rows (86 and 87) to the Thailand Thai League 1 2023-2024 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting pattern of an existing data row (row 2) onto the two
# new rows so that column A keeps the bold/bordered/centered style (s=1)
# and column E keeps the date-time number format (s=2), matching the rest
# of the table, without creating any duplicate style entries.
$ws.Range("A2:V2").Copy()
$ws.Range("A86:V86").PasteSpecial(-4122)
$ws.Range("A87:V87").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Row 86 ----
$ws.Cells.Item(86, 1).Value = 85
$ws.Cells.Item(86, 2).Value = "thailand"
$ws.Cells.Item(86, 3).Value = "thai-league-1"
$ws.Cells.Item(86, 4).Value = "2023-2024"
$ws.Cells.Item(86, 5).Value = 45261.54166666666
$ws.Cells.Item(86, 6).Value = "Ratchaburi"
$ws.Cells.Item(86, 7).Value = 3
$ws.Cells.Item(86, 8).Value = "Chiangrai Utd"
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 1.7
$ws.Cells.Item(86, 11).Value = "26/11/2023 13:12"
$ws.Cells.Item(86, 12).Value = 1.65
$ws.Cells.Item(86, 13).Value = "01/12/2023 12:55"
$ws.Cells.Item(86, 14).Value = 3.81
$ws.Cells.Item(86, 15).Value = "26/11/2023 13:12"
$ws.Cells.Item(86, 16).Value = 3.79
$ws.Cells.Item(86, 17).Value = "01/12/2023 12:59"
$ws.Cells.Item(86, 18).Value = 4.77
$ws.Cells.Item(86, 19).Value = "26/11/2023 13:12"
$ws.Cells.Item(86, 20).Value = 5.59
$ws.Cells.Item(86, 21).Value = "01/12/2023 12:59"
$ws.Cells.Item(86, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/ratchaburi-chiangrai-utd/fmqa26EH/"

# ---- Row 87 ----
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "thailand"
$ws.Cells.Item(87, 3).Value = "thai-league-1"
$ws.Cells.Item(87, 4).Value = "2023-2024"
$ws.Cells.Item(87, 5).Value = 45261.54166666666
$ws.Cells.Item(87, 6).Value = "Uthai Thani"
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = "Sukhothai"
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 1.85
$ws.Cells.Item(87, 11).Value = "26/11/2023 11:42"
$ws.Cells.Item(87, 12).Value = 2.24
$ws.Cells.Item(87, 13).Value = "01/12/2023 12:54"
$ws.Cells.Item(87, 14).Value = 3.79
$ws.Cells.Item(87, 15).Value = "26/11/2023 11:42"
$ws.Cells.Item(87, 16).Value = 3.62
$ws.Cells.Item(87, 17).Value = "01/12/2023 12:53"
$ws.Cells.Item(87, 18).Value = 3.95
$ws.Cells.Item(87, 19).Value = "26/11/2023 11:42"
$ws.Cells.Item(87, 20).Value = 3.11
$ws.Cells.Item(87, 21).Value = "01/12/2023 12:54"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/thailand/thai-league-1/uthai-thani-sukhothai/4bjn58Ub/"
